# The "brand_slug" column header (E1) is removed from the export;
# the header that used to sit in F1 ("brand_description") takes its place,
# and the now-trailing F1 cell is cleared out entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "brand_description"
$ws.Range("F1").ClearContents()
